{"js": "// Update the answer table in place: each of the 5 \"answer\" rows of the\n// worksheet (rows 0, 4, 8, 12, 16 \u2014 the other rows are blank work rows)\n// gets its 5 cell texts replaced with a freshly generated set of\n// two-digit-divided-by-one-digit problems/answers.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// New text for each of the 5 populated rows, keyed by the row's current\n// (pre-edit) first-cell text so this keeps working even if row order in\n// the table ever shifts.\nconst rowReplacements = {\n  \"51\u00f76=8, 3\": [\"76\u00f78=9, 4\", \"10\u00f78=1, 2\", \"51\u00f76=8, 3\", \"43\u00f72=21, 1\", \"34\u00f74=8, 2\"],\n  \"73\u00f76=12, 1\": [\"59\u00f73=19, 2\", \"43\u00f77=6, 1\", \"57\u00f74=14, 1\", \"45\u00f72=22, 1\", \"44\u00f79=4, 8\"],\n  \"88\u00f76=14, 4\": [\"89\u00f76=14, 5\", \"20\u00f79=2, 2\", \"49\u00f73=16, 1\", \"64\u00f74=16, 0\", \"56\u00f79=6, 2\"],\n  \"37\u00f76=6, 1\": [\"18\u00f76=3, 0\", \"90\u00f78=11, 2\", \"44\u00f76=7, 2\", \"53\u00f72=26, 1\", \"21\u00f73=7, 0\"],\n  \"11\u00f72=5, 1\": [\"44\u00f78=5, 4\", \"83\u00f79=9, 2\", \"71\u00f76=11, 5\", \"62\u00f72=31, 0\", \"67\u00f72=33, 1\"],\n};\n\nconst newValues = table.values.map((row) => {\n  const key = row[0];\n  return Object.prototype.hasOwnProperty.call(rowReplacements, key)\n    ? rowReplacements[key]\n    : row;\n});\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the answer table in place: each of the 5 \"answer\" rows of the\n# worksheet (rows 1, 5, 9, 13, 17 -- the other rows are blank work rows)\n# gets its 5 cell texts replaced with a freshly generated set of\n# two-digit-divided-by-one-digit problems/answers.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowValues = @{\n    1  = @(\"76\u00f78=9, 4\", \"10\u00f78=1, 2\", \"51\u00f76=8, 3\", \"43\u00f72=21, 1\", \"34\u00f74=8, 2\")\n    5  = @(\"59\u00f73=19, 2\", \"43\u00f77=6, 1\", \"57\u00f74=14, 1\", \"45\u00f72=22, 1\", \"44\u00f79=4, 8\")\n    9  = @(\"89\u00f76=14, 5\", \"20\u00f79=2, 2\", \"49\u00f73=16, 1\", \"64\u00f74=16, 0\", \"56\u00f79=6, 2\")\n    13 = @(\"18\u00f76=3, 0\", \"90\u00f78=11, 2\", \"44\u00f76=7, 2\", \"53\u00f72=26, 1\", \"21\u00f73=7, 0\")\n    17 = @(\"44\u00f78=5, 4\", \"83\u00f79=9, 2\", \"71\u00f76=11, 5\", \"62\u00f72=31, 0\", \"67\u00f72=33, 1\")\n}\n\nforeach ($rowIndex in $rowValues.Keys) {\n    $values = $rowValues[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
